$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the surviving rows (2-4) with the new "last cases" values
$ws.Range("A2").Value = 11
$ws.Range("B2").Value = 28

$ws.Range("A3").Value = 12
$ws.Range("B3").Value = 17

$ws.Range("A4").Value = 21
$ws.Range("B4").Value = 12

# Remove the now-unused rows 5-8 entirely
$ws.Rows("5:8").Delete() | Out-Null
